$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -87
$ws.Range("N2").ClearContents()
# Row 42
$ws.Range("H42").Value = 183
$ws.Range("I42").Value = 120
$ws.Range("J42").Value = 309
$ws.Range("K42").Value = 360
$ws.Range("L42").Value = 927
$ws.Range("M42").Value = -130
$ws.Range("N42").Value = -1387
# Row 51
$ws.Range("H51").Value = 4645.7617
$ws.Range("I51").Value = 3566.6667
$ws.Range("J51").Value = 4825.6113
$ws.Range("K51").Value = 3566.6667
$ws.Range("L51").Value = 4825.6113
$ws.Range("M51").Value = -3082.6667
$ws.Range("N51").Value = -5793.6113
# Row 138
$ws.Range("H138").Value = 6748686.5
$ws.Range("I138").Value = 1906886.5
$ws.Range("J138").Value = 10207115
$ws.Range("K138").Value = 5720659.5
$ws.Range("L138").Value = 30621345
$ws.Range("M138").Value = -5715519.5
$ws.Range("N138").Value = -30631625

$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 10745
$ws.Range("I28").Value = 4326.6665
$ws.Range("J28").Value = 30000
$ws.Range("K28").Value = 4326.6665
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = -4134.6665
$ws.Range("N28").Value = -30384
# Row 31
$ws.Range("H31").Value = 2284.4
$ws.Range("I31").Value = 2284.4
$ws.Range("K31").Value = 2284.4
$ws.Range("M31").Value = -1990.4
# Row 41
$ws.Range("H41").Value = 27031
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 27031
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 27031
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -27859
# Row 99
$ws.Range("H99").Value = 10745
$ws.Range("I99").Value = 4326.6665
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 4326.6665
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -1331.6665
$ws.Range("N99").Value = -35990
# Row 122
$ws.Range("H122").Value = 1792.125
$ws.Range("I122").Value = 1548.1428
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 4644.428400000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -2194.428400000001
$ws.Range("N122").Value = -15400

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 1152
$ws.Range("I8").Value = 304
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 304
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = -164
$ws.Range("N8").Value = -2280
# Row 64
$ws.Range("H64").Value = 897.8570999999999
$ws.Range("J64").Value = 435.58823
$ws.Range("L64").Value = 435.58823
$ws.Range("N64").Value = -885.5882300000001
# Row 67
$ws.Range("H67").Value = 897.8570999999999
$ws.Range("J67").Value = 435.58823
$ws.Range("L67").Value = 435.58823
$ws.Range("N67").Value = -1995.58823
# Row 122
$ws.Range("H122").Value = 32250
$ws.Range("J122").Value = 32250
$ws.Range("L122").Value = 32250
$ws.Range("N122").Value = -42050

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 16178.833
$ws.Range("I62").Value = 19137.072
$ws.Range("J62").Value = 5825
$ws.Range("K62").Value = 19137.072
$ws.Range("L62").Value = 5825
$ws.Range("M62").Value = -18513.072
$ws.Range("N62").Value = -7073
# Row 65
$ws.Range("H65").Value = 16178.833
$ws.Range("I65").Value = 19137.072
$ws.Range("J65").Value = 5825
$ws.Range("K65").Value = 95685.36
$ws.Range("L65").Value = 29125
$ws.Range("M65").Value = -92565.36
$ws.Range("N65").Value = -35365
# Row 99
$ws.Range("H99").Value = 20835840
$ws.Range("I99").Value = 2425.25
$ws.Range("J99").Value = 83336090
$ws.Range("K99").Value = 2425.25
$ws.Range("L99").Value = 83336090
$ws.Range("M99").Value = -927.25
$ws.Range("N99").Value = -83339086
# Row 126
$ws.Range("H126").Value = 20835840
$ws.Range("I126").Value = 2425.25
$ws.Range("J126").Value = 83336090
$ws.Range("K126").Value = 7275.75
$ws.Range("L126").Value = 250008270
$ws.Range("M126").Value = -4805.75
$ws.Range("N126").Value = -250013210

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 1663
$ws.Range("J75").Value = 1953.75
$ws.Range("L75").Value = 5861.25
$ws.Range("N75").Value = -7857.25
# Row 78
$ws.Range("H78").Value = 1663
$ws.Range("J78").Value = 1953.75
$ws.Range("L78").Value = 17583.75
$ws.Range("N78").Value = -27567.75
# Row 99
$ws.Range("H99").Value = 2177.2727
$ws.Range("I99").Value = 1658.3334
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 4975.0002
$ws.Range("L99").Value = 8400
$ws.Range("M99").Value = -2729.0002
$ws.Range("N99").Value = -12892
# Row 132
$ws.Range("H132").Value = 1778.2222
$ws.Range("J132").Value = 2040
$ws.Range("L132").Value = 18360
$ws.Range("N132").Value = -23420

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5533.5713
$ws.Range("I70").Value = 5581.552
$ws.Range("J70").Value = 5301.6665
$ws.Range("K70").Value = 5581.552
$ws.Range("L70").Value = 5301.6665
$ws.Range("M70").Value = -5311.552
$ws.Range("N70").Value = -5841.6665
# Row 73
$ws.Range("H73").Value = 5533.5713
$ws.Range("I73").Value = 5581.552
$ws.Range("J73").Value = 5301.6665
$ws.Range("K73").Value = 5581.552
$ws.Range("L73").Value = 5301.6665
$ws.Range("M73").Value = -4645.552
$ws.Range("N73").Value = -7173.6665
# Row 122
$ws.Range("H122").Value = 1837.6666
$ws.Range("I122").Value = 2053.5
$ws.Range("J122").Value = 1406
$ws.Range("K122").Value = 6160.5
$ws.Range("L122").Value = 4218
$ws.Range("M122").Value = -3710.5
$ws.Range("N122").Value = -9118
# Row 126
$ws.Range("H126").Value = 2948.05
$ws.Range("I126").Value = 2796.2
$ws.Range("K126").Value = 8388.599999999999
$ws.Range("M126").Value = -5918.599999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 13152.5
$ws.Range("I22").Value = 425
$ws.Range("J22").Value = 25880
$ws.Range("K22").Value = 425
$ws.Range("L22").Value = 25880
$ws.Range("M22").Value = -130
$ws.Range("N22").Value = -26470
# Row 27
$ws.Range("H27").Value = 13152.5
$ws.Range("I27").Value = 425
$ws.Range("J27").Value = 25880
$ws.Range("K27").Value = 425
$ws.Range("L27").Value = 25880
$ws.Range("M27").Value = -318
$ws.Range("N27").Value = -26094
# Row 132
$ws.Range("H132").Value = 5483.615
$ws.Range("I132").Value = 4051
$ws.Range("J132").Value = 6120.3335
$ws.Range("K132").Value = 12153
$ws.Range("L132").Value = 18361.0005
$ws.Range("M132").Value = -9623
$ws.Range("N132").Value = -23421.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
